$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")

# Insert a new column B ("Week_Start_Date") before the ASIN column,
# shifting ASIN..is_holiday_week from C..I to D..J.
$ws.Columns.Item(2).Insert()

# Force the new column to Text format so the date strings we write are
# kept as literal text instead of being auto-converted to date serials.
$ws.Columns.Item(2).NumberFormat = "@"

# Header for the new column.
$ws.Range("B1").Value = "Week_Start_Date"

# Week_Start_Date values for each week row (2 through 17).
$weekStartDates = @(
    "2025-01-05",
    "2025-01-12",
    "2025-01-19",
    "2025-01-26",
    "2025-02-02",
    "2025-02-09",
    "2025-02-16",
    "2025-02-23",
    "2025-03-02",
    "2025-03-09",
    "2025-03-16",
    "2025-03-23",
    "2025-03-30",
    "2025-04-06",
    "2025-04-13",
    "2025-04-20"
)

for ($i = 0; $i -lt $weekStartDates.Length; $i++) {
    $row = $i + 2
    $ws.Range("B$row").Value = $weekStartDates[$i]
}

# Strip the leading zero from the week labels for weeks 1-9 (W01..W09 -> W1..W9).
# Weeks 10-16 (rows 11-17) keep their existing labels (W10..W16) unchanged.
for ($week = 1; $week -le 9; $week++) {
    $row = $week + 1
    $ws.Range("A$row").Value = "W$week"
}

# Convert the is_holiday_week column (now column J) from numeric 0 to boolean FALSE.
for ($row = 2; $row -le 17; $row++) {
    $ws.Range("J$row").Value = $false
}
